$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "Registro de preço para futura e eventual contratação de empresa especializada na realização de exames patológicos junto a Secretaria de Saúde deste município"
$ws.Cells.Item(2, 2).Value = "147/2023"
$ws.Cells.Item(2, 3).Value = "063/2023"
$ws.Cells.Item(2, 4).Value = 45159
$ws.Cells.Item(2, 5).Value = 45174

# Row 3
$ws.Cells.Item(3, 1).Value = "Contratação de empresa com profissional médico na área de ortopedia, destinado à prestação de serviços junto a Secretaria Municipal de Saúde e Fundação Coronel João de Almeida, neste município"
$ws.Cells.Item(3, 2).Value = "190/2023"
$ws.Cells.Item(3, 3).Value = "081/2023"
$ws.Cells.Item(3, 4).Value = 45230
$ws.Cells.Item(3, 5).Value = 45250

# Row 4
$ws.Cells.Item(4, 1).Value = "contratação de empresa com profissional cirurgião dentista para prestar serviços em unidades básicas de saúde"
$ws.Cells.Item(4, 2).Value = "161/2023"
$ws.Cells.Item(4, 3).Value = "069/2023"
$ws.Cells.Item(4, 4).Value = 45184
$ws.Cells.Item(4, 5).Value = 45205

# Row 5
$ws.Cells.Item(5, 1).Value = "Contratação de empresa com profissional médico especializado na área de cardiologia, para realização de consultas em unidade de saúde deste município e realização de procedimentos em consultório clínico da contratada"
$ws.Cells.Item(5, 2).Value = "139/2023"
$ws.Cells.Item(5, 3).Value = "061/2023"
$ws.Cells.Item(5, 4).Value = "15/08/2023"
$ws.Cells.Item(5, 5).Value = 45167

# Apply the date number format once on a "master" cell, then propagate it
# via copy/paste-special (format only) so every date cell shares a single
# style record instead of each write minting its own.
$ws.Cells.Item(2, 4).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(2, 4).Copy()
$ws.Range("D2:E4").PasteSpecial(-4122)
$ws.Cells.Item(5, 5).PasteSpecial(-4122)

# D5 keeps a text-formatted date (matching the existing text style used by
# D1/E1), so copy that style over instead of minting a new one.
$ws.Cells.Item(1, 4).Copy()
$ws.Cells.Item(5, 4).PasteSpecial(-4122)

$excel.CutCopyMode = $false

$ws.Range("E5").Select()
